$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 1).Value = "QSR"
$ws.Cells.Item(6, 2).Value = "Quasar"

$ws.Range("C8").Select()
